# Actualiza la base de datos del Estado de Cuenta:
# Se elimina el agrupamiento por periodo (todos 1809 primero, luego todos 1810)
# y se reordena para agrupar por trabajador, mostrando primero el periodo 1810
# (el mas reciente) seguido del 1809 para cada trabajador.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Nuevos valores (Documento, Nombre, Periodo) para las filas 16 a 21
$datos = @(
    @("73164755",    "CAMILO YEPES CADENA",            "1810"),
    @("73164755",    "CAMILO YEPES CADENA",            "1809"),
    @("1143354313",  "WILMER JAVIER DOMINGUEZ CASTRO",  "1810"),
    @("1143354313",  "WILMER JAVIER DOMINGUEZ CASTRO",  "1809"),
    @("1002250621",  "ANDREINA PAOLA ROMERO NARVAEZ",   "1810"),
    @("1002250621",  "ANDREINA PAOLA ROMERO NARVAEZ",   "1809")
)

$fila = 16
foreach ($registro in $datos) {
    $ws.Cells.Item($fila, 3).Value = $registro[0]
    $ws.Cells.Item($fila, 4).Value = $registro[1]
    $ws.Cells.Item($fila, 5).Value = $registro[2]
    $fila = $fila + 1
}
